# Apply updated odds/values from the 2024-10-12 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5

# Row 8
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.9

# Row 9
$ws.Range("G9").Value = 2.3
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 2.88
$ws.Range("L9").Value = 3.4
$ws.Range("AI9").Value = 11
$ws.Range("AU9").Value = 7.5
$ws.Range("BD9").Value = 126

# Row 11
$ws.Range("G11").Value = 3.65
$ws.Range("H11").Value = 3.05
$ws.Range("I11").Value = 2.05
$ws.Range("J11").Value = 4.15
$ws.Range("L11").Value = 2.67
$ws.Range("N11").Value = 6.5
$ws.Range("Q11").Value = 2.22
$ws.Range("R11").Value = 1.52
$ws.Range("S11").Value = 1.45
$ws.Range("T11").Value = 2.37
$ws.Range("W11").Value = 8.75
$ws.Range("X11").Value = 18.5
$ws.Range("Y11").Value = 13
$ws.Range("Z11").Value = 55
$ws.Range("AB11").Value = 50
$ws.Range("AC11").Value = 7.2
$ws.Range("AD11").Value = 6
$ws.Range("AE11").Value = 17
$ws.Range("AH11").Value = 8.75
$ws.Range("AJ11").Value = 18.5
$ws.Range("AK11").Value = 19.5
$ws.Range("AN11").Value = 5.3
$ws.Range("AO11").Value = 21
$ws.Range("AP11").Value = 29
$ws.Range("AU11").Value = 7.5
$ws.Range("AV11").Value = 75
$ws.Range("AW11").Value = 3.75
$ws.Range("AX11").Value = 10.75
$ws.Range("AZ11").Value = 45
$ws.Range("BB11").Value = 300

# Row 15
$ws.Range("M15").Value = 1.03
$ws.Range("O15").Value = 1.27

# Row 16
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 9
$ws.Range("O16").Value = 1.33

# Row 17
$ws.Range("M17").Value = 1.04
$ws.Range("O17").Value = 1.3

# Row 24
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10

# Row 29
$ws.Range("O29").Value = 1.24
$ws.Range("R29").Value = 2.02
$ws.Range("AB29").Value = 26
$ws.Range("AF29").Value = 50
$ws.Range("AU29").Value = 6.9

# Row 30
$ws.Range("G30").Value = 1.53
$ws.Range("H30").Value = 4.3
$ws.Range("I30").Value = 4.6
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 2.57
$ws.Range("L30").Value = 4.6
$ws.Range("N30").Value = 9.75
$ws.Range("R30").Value = 2.65
$ws.Range("S30").Value = 1.24
$ws.Range("U30").Value = 1.52
$ws.Range("V30").Value = 2.37
$ws.Range("W30").Value = 11.25
$ws.Range("X30").Value = 10
$ws.Range("Z30").Value = 13
$ws.Range("AA30").Value = 11
$ws.Range("AB30").Value = 17.5
$ws.Range("AC30").Value = 9.75
$ws.Range("AD30").Value = 9.5
$ws.Range("AE30").Value = 13.5
$ws.Range("AF30").Value = 40
$ws.Range("AG30").Value = 20
$ws.Range("AH30").Value = 35
$ws.Range("AI30").Value = 15.5
$ws.Range("AJ30").Value = 80
$ws.Range("AK30").Value = 37
$ws.Range("AL30").Value = 32
$ws.Range("AM30").Value = 200
$ws.Range("AN30").Value = 3.9
$ws.Range("AO30").Value = 7.2
$ws.Range("AQ30").Value = 19
$ws.Range("AR30").Value = 35
$ws.Range("AU30").Value = 6.8
$ws.Range("AW30").Value = 7
$ws.Range("AX30").Value = 24
$ws.Range("AY30").Value = 23

# Row 31
$ws.Range("G31").Value = 1.65
$ws.Range("H31").Value = 3.65
$ws.Range("I31").Value = 4.45
$ws.Range("J31").Value = 2.25
$ws.Range("L31").Value = 4.7
$ws.Range("N31").Value = 7.8
$ws.Range("O31").Value = 1.26
$ws.Range("P31").Value = 3.5
$ws.Range("Q31").Value = 1.8
$ws.Range("R31").Value = 1.95
$ws.Range("X31").Value = 8
$ws.Range("Y31").Value = 8
$ws.Range("Z31").Value = 13
$ws.Range("AA31").Value = 13
$ws.Range("AC31").Value = 7.8
$ws.Range("AD31").Value = 7.2
$ws.Range("AG31").Value = 13.5
$ws.Range("AH31").Value = 27
$ws.Range("AI31").Value = 14.5
$ws.Range("AJ31").Value = 75
$ws.Range("AK31").Value = 40
$ws.Range("AL31").Value = 45
$ws.Range("AN31").Value = 3.55
$ws.Range("AO31").Value = 8.25
$ws.Range("AP31").Value = 17.5
$ws.Range("AQ31").Value = 28
$ws.Range("AW31").Value = 6.3
$ws.Range("AX31").Value = 25
$ws.Range("AY31").Value = 30
$ws.Range("AZ31").Value = 150
$ws.Range("BA31").Value = 175
$ws.Range("BB31").Value = 400

# Row 32
$ws.Range("G32").Value = 2.5
$ws.Range("H32").Value = 3.15
$ws.Range("J32").Value = 3.15
$ws.Range("K32").Value = 2.07
$ws.Range("L32").Value = 3.2
$ws.Range("N32").Value = 7.3
$ws.Range("W32").Value = 8.5
$ws.Range("Z32").Value = 28
$ws.Range("AA32").Value = 21
$ws.Range("AB32").Value = 28
$ws.Range("AC32").Value = 7.3
$ws.Range("AD32").Value = 6.2
$ws.Range("AE32").Value = 12.5
$ws.Range("AG32").Value = 9.5
$ws.Range("AH32").Value = 14.5
$ws.Range("AJ32").Value = 32
$ws.Range("AK32").Value = 20
$ws.Range("AL32").Value = 27
$ws.Range("AN32").Value = 4.55
$ws.Range("AO32").Value = 14
$ws.Range("AQ32").Value = 60
$ws.Range("AR32").Value = 100
$ws.Range("AU32").Value = 6.7
$ws.Range("AV32").Value = 55
$ws.Range("AX32").Value = 14
$ws.Range("AY32").Value = 20
$ws.Range("AZ32").Value = 60
$ws.Range("BA32").Value = 90
